{"js": "// Replace the multiplication problems in the practice table with the new\n// set of problems, while preserving the existing run/paragraph formatting\n// (font, size, justification) of every cell.\n//\n// The table has 20 rows x 5 columns; only 5 of the rows (0, 4, 9, 14, 19)\n// actually contain problem text - the rows in between are left blank for\n// students to write their answers. Every problem cell's text is replaced\n// in row-major, left-to-right order (matching the order the problems\n// appear in the document).\n\nconst newProblemsByRow = {\n  0: [\"337\u00d74=\", \"991\u00d79=\", \"430\u00d79=\", \"378\u00d75=\", \"946\u00d78=\"],\n  4: [\"893\u00d78=\", \"691\u00d78=\", \"409\u00d77=\", \"726\u00d74=\", \"825\u00d74=\"],\n  9: [\"301\u00d76=\", \"170\u00d78=\", \"894\u00d79=\", \"980\u00d74=\", \"525\u00d79=\"],\n  14: [\"435\u00d78=\", \"284\u00d74=\", \"173\u00d73=\", \"147\u00d74=\", \"449\u00d75=\"],\n  19: [\"387\u00d77=\", \"923\u00d75=\", \"274\u00d76=\", \"516\u00d76=\", \"244\u00d76=\"],\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const rowIndexStr of Object.keys(newProblemsByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = newProblemsByRow[rowIndex];\n\n  for (let colIndex = 0; colIndex < newValues.length; colIndex++) {\n    const cell = table.getCell(rowIndex, colIndex);\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    // Replace the paragraph's range text in place so the existing run\n    // formatting (font/size) and paragraph formatting (justification)\n    // carry over to the new text, instead of being reset.\n    const range = paragraphs.items[0].getRange();\n    range.insertText(newValues[colIndex], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the multiplication problems in the practice table with the new\n# set of problems, while preserving the existing run/paragraph formatting\n# (font, size, justification) of every cell.\n#\n# The table has 20 rows x 5 columns; only 5 of the rows (1, 5, 10, 15, 20 -\n# 1-based, as COM table rows/cells are 1-indexed) actually contain problem\n# text - the rows in between are left blank for students to write their\n# answers. Setting Cell.Range.Text keeps the existing run/paragraph\n# formatting of the single run already in each cell.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$newProblemsByRow = @{\n    1  = @(\"337\u00d74=\", \"991\u00d79=\", \"430\u00d79=\", \"378\u00d75=\", \"946\u00d78=\")\n    5  = @(\"893\u00d78=\", \"691\u00d78=\", \"409\u00d77=\", \"726\u00d74=\", \"825\u00d74=\")\n    10 = @(\"301\u00d76=\", \"170\u00d78=\", \"894\u00d79=\", \"980\u00d74=\", \"525\u00d79=\")\n    15 = @(\"435\u00d78=\", \"284\u00d74=\", \"173\u00d73=\", \"147\u00d74=\", \"449\u00d75=\")\n    20 = @(\"387\u00d77=\", \"923\u00d75=\", \"274\u00d76=\", \"516\u00d76=\", \"244\u00d76=\")\n}\n\nforeach ($rowIndex in $newProblemsByRow.Keys) {\n    $newValues = $newProblemsByRow[$rowIndex]\n    for ($colIndex = 1; $colIndex -le $newValues.Length; $colIndex++) {\n        $cell = $table.Cell($rowIndex, $colIndex)\n        $cell.Range.Text = $newValues[$colIndex - 1]\n    }\n}\n"}
